# Cnc-Calculators-V.2/Database.xlsx
# "we are now dealing with non number inputs, they should never make their
#  way to the excel file" -- append rows that exercise non-numeric / blank
# input handling on every calculator sheet, plus move the selection on the
# last (active) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Cutting Speed"  (B2:E5 -> B2:E10)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cutting Speed")
$ws.Range("B6").Value = "Hei"
$ws.Range("C6").Value = "Ærling"
$ws.Range("D6").Value = "the "
$ws.Range("E6").Value = "Lærling"

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 80
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 0.12

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

# ---------------------------------------------------------------------
# Sheet 2: "Material Removal Rate"  (B2:F4 -> B2:F25)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Material Removal Rate")

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = 602
$ws.Range("E5").Value = 36.12
$ws.Range("F5").Value = "cm³/min"

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 602
$ws.Range("E6").Value = 72.23999999999999
$ws.Range("F6").Value = "cm³/min"

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 602
$ws.Range("E7").Value = 36.12
$ws.Range("F7").Value = "cm³/min"

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 60
$ws.Range("D8").Value = 602
$ws.Range("E8").Value = 72.23999999999999
$ws.Range("F8").Value = "cm³/min"

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 602
$ws.Range("E9").Value = 36.12
$ws.Range("F9").Value = "cm³/min"

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 60
$ws.Range("D10").Value = 1208
$ws.Range("E10").Value = 72.48
$ws.Range("F10").Value = "cm³/min"

$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 60
$ws.Range("D11").Value = 602
$ws.Range("E11").Value = 36.12
$ws.Range("F11").Value = "cm³/min"

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = 301
$ws.Range("E12").Value = 18.06
$ws.Range("F12").Value = "cm³/min"

$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = 301
$ws.Range("E13").Value = 36.12
$ws.Range("F13").Value = "cm³/min"

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 60
$ws.Range("D14").Value = 301
$ws.Range("E14").Value = 18.06
$ws.Range("F14").Value = "cm³/min"

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 60
$ws.Range("D15").Value = 602
$ws.Range("E15").Value = 36.12
$ws.Range("F15").Value = "cm³/min"

$ws.Range("B16").Value = "ehi"
$ws.Range("C16").Value = 60
$ws.Range("D16").Value = 602
$ws.Range("E16").Value = "Please input values"
$ws.Range("F16").Value = "cm³/min"

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 602
$ws.Range("E17").Value = 36.12
$ws.Range("F17").Value = "cm³/min"

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "Please input values"
$ws.Range("F18").Value = "cm³/min"

$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = "cm³/min"

$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = "cm³/min"

$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = "cm³/min"

$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = "cm³/min"

$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 60
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = "cm³/min"

$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 602
$ws.Range("E24").Value = 36.12
$ws.Range("F24").Value = "cm³/min"

$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 60
$ws.Range("D25").Value = 602
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = "cm³/min"

# ---------------------------------------------------------------------
# Sheet 3: "Helix Angle"  (B2:F5 -> B2:F11)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Helix Angle")

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "Please input values"
$ws.Range("F6").Value = "°"

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "°"

$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "°"

$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "°"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 0.06
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "°"

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 0.06
$ws.Range("E11").Value = 1.09
$ws.Range("F11").Value = "°"

# ---------------------------------------------------------------------
# Sheet 4: "Ramp Angle"  (B2:E6 -> B2:E10)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ramp Angle")

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "Please input values"
$ws.Range("E7").Value = "°"

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "°"

$ws.Range("B9").Value = 1000
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "°"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "°"

# ---------------------------------------------------------------------
# Sheet 5: "Surface Roughness"  (B2:D9 -> B2:D12)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Surface Roughness")

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

$ws.Range("B11").Value = 0.25
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 0.25
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0

# This sheet is the active/visible one -- match the saved selection.
$ws.Activate() | Out-Null
$ws.Range("C12").Select() | Out-Null
